# Update version string across the workbook for new release:
# "mines - January 30 (built on February 02 2026 12.49.33 EST)"
#   -> "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# About sheet: A2 "Version: ..."
$wsAbout.Range("A2").Value = "Version: " + $newVersion

# About sheet: A6 Recommended Citation text referencing the version
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Ashton Coal Mine, Australia, M0007, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Boundaries and methane sources sheet: S2:S17 build_version column values
for ($r = 2; $r -le 17; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
